$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the phone/name sample value in A2
$ws.Range("A2").Value = "0993000001"

# Clear SHORT_DESC (B2) and DESC (C2) sample values
$ws.Range("B2").ClearContents()
$ws.Range("C2").ClearContents()

# Set SKU (D2) to the same text value, preserving leading zero as text
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0993000001"

# Update PRICE (E2)
$ws.Range("E2").Value = 50000

# Set CATEGORY_ID (F2)
$ws.Range("F2").Value = 3

# Update BRAND (G2)
$ws.Range("G2").Value = "VNP"

# Update selection to C2
$ws.Range("C2").Select()
